# Reflektor Presentation Deck - minor typo fix
#
# 1) Bump the cached "Date and time" footer field (datetimeFigureOut,
#    fld id {1D8BD707-D9CF-40AE-B4C6-C98DA3205C09}) from 05/03/2025 to
#    05/04/2025 on the Slide Master and on every Slide Layout.
# 2) On slide 5 ("Coding Approach"), widen the title textbox and
#    upper-case its text to "CODING APPROACH".

$p = $ppt.ActivePresentation

$oldDate = "05/03/2025"
$newDate = "05/04/2025"

# --- Slide Master date placeholder ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every Slide Layout's date placeholder ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 5: "Coding Approach" title textbox ---
$slide = $p.Slides.Item(5)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $shp = $slide.Shapes.Item($si)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "Coding Approach") {
            # Set the text first (autosized height is computed against the
            # shape's current width), then widen the box.
            $shp.TextFrame.TextRange.Text = "CODING APPROACH"
            $shp.Width = 629.703779527559
        }
    }
}
